$d = $word.ActiveDocument

# Locate the paragraph holding the "{kepala_desa}" placeholder - the
# "Pj. Kepala Desa {nama_desa}" signature line right above the blank
# signature-line paragraphs. The fix adds a missing "{nip}" line
# (right-aligned under the signature via a run of tabs) directly below it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*{kepala_desa}*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph containing {kepala_desa}"
}

# Formatting shared by every run in the paragraph (matches the
# surrounding "{kepala_desa}" line: Times New Roman, 12pt).
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr>'
$pPr = '<w:pPr><w:spacing w:after="120"/><w:ind w:left="288"/>' + $rPr + '</w:pPr>'

# Build 8 separate "just a tab" runs, then a final run that carries the
# 9th tab plus the "{nip}" text - this mirrors how the preceding
# "{kepala_desa}" line is itself built up out of one-tab-per-run runs.
$tabRun = '<w:r>' + $rPr + '<w:tab/></w:r>'
$runsXml = ($tabRun * 8) + '<w:r>' + $rPr + '<w:tab/><w:t>{nip}</w:t></w:r>'

$paraXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $pPr + $runsXml + '</w:p>'

# Insert a brand-new empty paragraph right after the "{kepala_desa}" line,
# then replace its (empty) contents with the XML built above so the tabs
# land in the document as real <w:tab/> run children rather than literal
# tab characters inside a <w:t>.
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Range.InsertXML($paraXml)
